$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3103.8333
$ws.Range("I32").Value = 1481.8572
$ws.Range("J32").Value = 4136
$ws.Range("K32").Value = 1481.8572
$ws.Range("L32").Value = 4136
$ws.Range("M32").Value = -1155.8572
$ws.Range("N32").Value = -4788

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3604.4443
$ws.Range("I106").Value = 1639.5652
$ws.Range("J106").Value = 14902.5
$ws.Range("K106").Value = 1639.5652
$ws.Range("L106").Value = 14902.5
$ws.Range("M106").Value = -1008.5652
$ws.Range("N106").Value = -16164.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 26628.846
$ws.Range("I111").Value = 11897.777
$ws.Range("J111").Value = 59773.75
$ws.Range("K111").Value = 35693.331
$ws.Range("L111").Value = 179321.25
$ws.Range("M111").Value = -32626.331
$ws.Range("N111").Value = -185455.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3641.3137
$ws.Range("I132").Value = 3066.8542
$ws.Range("K132").Value = 9200.562600000001
$ws.Range("M132").Value = -6670.562600000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 15627374
$ws.Range("I137").Value = 52633590
$ws.Range("K137").Value = 157900770
$ws.Range("M137").Value = -157898220

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2755.8965
$ws.Range("J138").Value = 3367.8
$ws.Range("L138").Value = 10103.4
$ws.Range("N138").Value = -20383.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 372.04
$ws.Range("I32").Value = 353.91666
$ws.Range("K32").Value = 353.91666
$ws.Range("M32").Value = -66.91665999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1672.9166
$ws.Range("I74").Value = 1470.5714
$ws.Range("K74").Value = 1470.5714
$ws.Range("M74").Value = -596.5714

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1672.9166
$ws.Range("I77").Value = 1470.5714
$ws.Range("K77").Value = 7352.857
$ws.Range("M77").Value = -2984.857

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 15290.071
$ws.Range("J102").Value = 17399.6
$ws.Range("L102").Value = 17399.6
$ws.Range("N102").Value = -20643.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 4151.3335
$ws.Range("I122").Value = 3571.1875
$ws.Range("K122").Value = 10713.5625
$ws.Range("M122").Value = -8263.5625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 55556950
$ws.Range("I80").Value = 897.6667
$ws.Range("K80").Value = 897.6667
$ws.Range("M80").Value = 100.3333

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H83").Value = 55556950
$ws.Range("I83").Value = 897.6667
$ws.Range("K83").Value = 4488.3335
$ws.Range("M83").Value = 503.6665000000003

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 71434216
$ws.Range("I86").Value = 100006100
$ws.Range("K86").Value = 100006100
$ws.Range("M86").Value = -100004977

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 71434216
$ws.Range("I89").Value = 100006100
$ws.Range("K89").Value = 500030500
$ws.Range("M89").Value = -500024884

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1942.5952
$ws.Range("I94").Value = 1081.25
$ws.Range("J94").Value = 3091.0557
$ws.Range("K94").Value = 1081.25
$ws.Range("L94").Value = 3091.0557
$ws.Range("M94").Value = -630.25
$ws.Range("N94").Value = -3993.0557

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2273.5715
$ws.Range("I134").Value = 1871.5385
$ws.Range("J134").Value = 7500
$ws.Range("K134").Value = 5614.6155
$ws.Range("L134").Value = 22500
$ws.Range("M134").Value = -3079.6155
$ws.Range("N134").Value = -27570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H41").Value = 9026.5
$ws.Range("I41").Value = 3483.4546
$ws.Range("K41").Value = 3483.4546
$ws.Range("M41").Value = -3055.4546

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 313.83334
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 1012.375
$ws.Range("J50").Value = 1012.375
$ws.Range("L50").Value = 3037.125
$ws.Range("N50").Value = -3999.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H53").Value = 1012.375
$ws.Range("J53").Value = 1012.375
$ws.Range("L53").Value = 3037.125
$ws.Range("N53").Value = -3999.125

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 494.07693
$ws.Range("I92").Value = 192.4
$ws.Range("K92").Value = 577.2
$ws.Range("M92").Value = 670.8

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H114").Value = 1283.0667
$ws.Range("I114").Value = 1096.3334
$ws.Range("J114").Value = 1563.1666
$ws.Range("K114").Value = 3289.0002
$ws.Range("L114").Value = 4689.4998
$ws.Range("M114").Value = -35.00019999999995
$ws.Range("N114").Value = -11197.4998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 2497.3333
$ws.Range("I134").Value = 1996.8
$ws.Range("J134").Value = 5000
$ws.Range("K134").Value = 5990.4
$ws.Range("L134").Value = 15000
$ws.Range("M134").Value = -920.3999999999996
$ws.Range("N134").Value = -25140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 313.83334
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 1192.1428
$ws.Range("I139").Value = 890.8333
$ws.Range("K139").Value = 2672.4999
$ws.Range("M139").Value = 2467.5001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 3006
$ws.Range("I140").Value = 2341.3333
$ws.Range("J140").Value = 5000
$ws.Range("K140").Value = 7023.999899999999
$ws.Range("L140").Value = 15000
$ws.Range("M140").Value = -1843.999899999999
$ws.Range("N140").Value = -25360

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 41000
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 41000
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 41000
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -42640

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4227.857
$ws.Range("I113").Value = 2996.6667
$ws.Range("J113").Value = 4563.636
$ws.Range("K113").Value = 2996.6667
$ws.Range("L113").Value = 4563.636
$ws.Range("M113").Value = -826.6667000000002
$ws.Range("N113").Value = -8903.636

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2187.1365
$ws.Range("I122").Value = 1900.9474
$ws.Range("K122").Value = 5702.8422
$ws.Range("M122").Value = -3252.8422

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H123").Value = 29931.666
$ws.Range("J123").Value = 29931.666
$ws.Range("L123").Value = 29931.666
$ws.Range("N123").Value = -34831.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8067765.5
$ws.Range("I22").Value = 3999
$ws.Range("J22").Value = 16131532
$ws.Range("K22").Value = 3999
$ws.Range("L22").Value = 16131532
$ws.Range("M22").Value = -3704
$ws.Range("N22").Value = -16132122

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 8067765.5
$ws.Range("I27").Value = 3999
$ws.Range("J27").Value = 16131532
$ws.Range("K27").Value = 3999
$ws.Range("L27").Value = 16131532
$ws.Range("M27").Value = -3892
$ws.Range("N27").Value = -16131746

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3213.7666
$ws.Range("I40").Value = 2859.7273
$ws.Range("J40").Value = 4187.375
$ws.Range("K40").Value = 2859.7273
$ws.Range("L40").Value = 4187.375
$ws.Range("M40").Value = -2723.7273
$ws.Range("N40").Value = -4459.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 17748.75
$ws.Range("I100").Value = 28497.5
$ws.Range("K100").Value = 28497.5
$ws.Range("M100").Value = -27956.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H119").Value = 0
$ws.Range("I119").Value = 0
$ws.Range("K119").Value = 0
$ws.Range("M119").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3763.879
$ws.Range("I122").Value = 2992.7058
$ws.Range("K122").Value = 8978.117400000001
$ws.Range("M122").Value = -6528.117400000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 6153.269
$ws.Range("I132").Value = 3606.7144
$ws.Range("K132").Value = 10820.1432
$ws.Range("M132").Value = -8290.143199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5888.278
$ws.Range("I136").Value = 1620
$ws.Range("K136").Value = 4860
$ws.Range("M136").Value = -2310

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H23").Value = 3426.4546
$ws.Range("I23").Value = 2083.625
$ws.Range("J23").Value = 7007.3335
$ws.Range("K23").Value = 2083.625
$ws.Range("L23").Value = 7007.3335
$ws.Range("M23").Value = -1854.625
$ws.Range("N23").Value = -7465.3335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 32000
$ws.Range("J40").Value = 32000
$ws.Range("L40").Value = 32000
$ws.Range("N40").Value = -32298

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 18949
$ws.Range("I43").Value = 18949
$ws.Range("K43").Value = 18949
$ws.Range("M43").Value = -18800

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 488.2857
$ws.Range("J107").Value = 501.5
$ws.Range("L107").Value = 1504.5
$ws.Range("N107").Value = -5344.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2340.2104
$ws.Range("I122").Value = 1748.0667
$ws.Range("K122").Value = 5244.2001
$ws.Range("M122").Value = -2794.2001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 109995.5
$ws.Range("J133").Value = 109995.5
$ws.Range("L133").Value = 109995.5
$ws.Range("N133").Value = -120115.5
